# Commit: "added open , unzip, clearing when close"
#
# The user opened the workbook, switched to the (previously empty) "Sheet1"
# worksheet, and typed two new values into it: "sdfz" into D7 and "zsd" into
# C2. (The shared-string table shows "sdfz" registered at the lower index,
# so it must have been entered first, even though D7 sorts after C2 in the
# sheet.) The selection was then left on B14 when the workbook was saved,
# which is also what makes "Sheet1" the active tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# Entry order matters: it controls shared-string allocation order.
$ws.Range("D7").Value = "sdfz"
$ws.Range("C2").Value = "zsd"

# Final selection left on B14 when the sheet/workbook was saved.
$ws.Range("B14").Select() | Out-Null
